$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.963.88'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.933.32'
$ws.Range("E3").Value = '  +3.57%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.622'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.36'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0877'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.137'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.10'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = '3.395.04'
$ws.Range("E14").Value = '  +3.66%  '
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").Value = '2.934.04'
$ws.Range("E16").Value = '  +3.61%  '
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").Value = '52.025.21'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.97%  '
$ws.Range("D22").Value = '0.0₃0980'
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("E26").Value = '  +11.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.38%  '
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  +16.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.106'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +17.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("E32").Value = '  -0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '37.08'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '52.91'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.67'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.80%  '
$ws.Range("E40").Value = '  +1.64%  '
$ws.Range("E41").Value = '  +5.09%  '
$ws.Range("E42").Value = '  +1.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.58%  '
$ws.Range("E45").Value = '  +1.89%  '
$ws.Range("D46").Value = '2.191.41'
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '110.93'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.248'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0350'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.948'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.93%  '
